$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1 / rId1) — update column F ("想去人数") values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1469
$ws1.Range("F7").Value = 139
$ws1.Range("F8").Value = 6179
$ws1.Range("F10").Value = 402
$ws1.Range("F11").Value = 113
$ws1.Range("F12").Value = 5056
$ws1.Range("F17").Value = 358
$ws1.Range("F18").Value = 62
$ws1.Range("F20").Value = 291
$ws1.Range("F22").Value = 3578
$ws1.Range("F23").Value = 147

# Sheet "全部类型" (sheet4 / rId4) — same underlying events, mirrored rows
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 1469
$ws4.Range("F8").Value = 139
$ws4.Range("F9").Value = 6179
$ws4.Range("F11").Value = 402
$ws4.Range("F12").Value = 113
$ws4.Range("F13").Value = 5056
$ws4.Range("F18").Value = 358
$ws4.Range("F19").Value = 62
$ws4.Range("F21").Value = 291
$ws4.Range("F23").Value = 3578
$ws4.Range("F25").Value = 147
